# Styling of requirement documents are updated
#
# Both edits add "their first name, last name," into the two "register a
# [guest/new] user" bullets on the User Requirements slides.

$p = $ppt.ActivePresentation

# --- Slide 4: "Functional Requirements" -> Account Requirements -------
# "A guest user shall be able to register to the system with a driver's
# license, ..." becomes "...with their first name, last name, a driver's
# license, ...". In the source this edit was typed as three runs: the
# unchanged lead-in, a single space that picked up a tr-TR language tag,
# and the unchanged rest of the sentence (now prefixed with the new
# clause). Rebuild the paragraph with InsertAfter so the run boundaries
# land in the same place.
$slide4 = $p.Slides.Item(4)
$guestShape = $slide4.Shapes.Item("Content Placeholder 2")
$guestRange = $guestShape.TextFrame.TextRange
$guestPara = $guestRange.Paragraphs(3, 1)

$guestPara.Text = "A guest user shall be able to register to the system with"
$guestSpace = $guestPara.InsertAfter(" ")
$guestSpace.LanguageID = "tr-TR"
$guestPara.InsertAfter("their first name, last name, a driver’s license, a bank card number, a phone number, an e-mail address and a valid password to become a member.") | Out-Null

# --- Slide 5: "Functional Requirements" -> Employee Requirements ------
# "An employee shall be able to register a new user with a driver's
# license, ..." becomes "...with their first name, last name, a driver's
# license, ...". This one stays a single run in the final XML, so stage
# the edit through a throwaway placeholder first -- assigning the final
# sentence directly would get diffed against the old run and re-split
# into multiple runs around the inserted words.
$slide5 = $p.Slides.Item(5)
$employeeShape = $slide5.Shapes.Item("Content Placeholder 2")
$employeeRange = $employeeShape.TextFrame.TextRange
$employeePara = $employeeRange.Paragraphs(9, 1)

$employeePara.Text = "placeholder"
$employeePara.Text = "An employee shall be able to register a new user with their first name, last name, a driver’s license, a credit card number, a phone number, an e-mail address and a valid password."
